# Add an "additional_brcs" slot to the Dataset class/sheet.
#
# This mirrors the xml diff:
#   - a new column is inserted at F ("additional_brcs"), shifting every
#     column from the old "repository" (F) onward one position to the right
#   - the sheet dimension grows from A1:R1 to A1:S1
#   - a new list-style data validation (same BRC list as column E) is
#     applied to the new F column
#   - the existing "datasetType" validation (previously on L) now targets M

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# Insert a new column before the current "repository" column (F),
# shifting repository..dataset_url one column to the right.
$ws.Columns.Item(6).Insert()

# Label the newly inserted column.
$ws.Range("F1").Value = "additional_brcs"

# Drop any validations Excel may have auto-extended onto the shifted
# columns so we can (re)create them in the same order as the target file.
$ws.Range("E2:E1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("M2:M1048576").Validation.Delete()

$xlValidateList = [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList
$xlValidAlertStop = [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop

# brc (E): unchanged, re-applied so it keeps its original position/order.
$brcRange = $ws.Range("E2:E1048576")
$brcRange.Validation.Add($xlValidateList, $xlValidAlertStop, 1, """CABBI,CBI,GLBRC,JBEI""")
$brcRange.Validation.ShowInput = $false
$brcRange.Validation.ShowError = $false

# additional_brcs (F): new validation, same dropdown list as brc.
$additionalBrcRange = $ws.Range("F2:F1048576")
$additionalBrcRange.Validation.Add($xlValidateList, $xlValidAlertStop, 1, """CABBI,CBI,GLBRC,JBEI""")
$additionalBrcRange.Validation.ShowInput = $false
$additionalBrcRange.Validation.ShowError = $false

# datasetType validation, now on M after the column shift.
$datasetTypeRange = $ws.Range("M2:M1048576")
$datasetTypeRange.Validation.Add($xlValidateList, $xlValidAlertStop, 1, """AS,GD,IM,ND,IP,FP,SM,MM,I""")
$datasetTypeRange.Validation.ShowInput = $false
$datasetTypeRange.Validation.ShowError = $false
